$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$headers = @("Company", "Role", "Location", "Platform", "Status", "Note", "Date", "Time", "Link")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# New data rows (Company, Role, Location, Platform, Status, Note, Date, Time, Link)
$rows = @(
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Manual Review", "Timeout 30000ms exceeded.", "19-01-26", "18:00", "https://www.linkedin.com/jobs/search/?currentJobId=4364382406&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Manual Review", "Timeout 30000ms exceeded.", "19-01-26", "18:01", "https://www.linkedin.com/jobs/search/?currentJobId=4364203279&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Manual Review", "Timeout 30000ms exceeded.", "19-01-26", "18:05", "https://www.linkedin.com/jobs/search/?currentJobId=4364382406&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland"),
    @("Test Company Ltd", "Junior Software Engineer", "Ireland", "LinkedIn", "Manual Review", "Timeout 30000ms exceeded.", "19-01-26", "18:06", "https://www.linkedin.com/jobs/search/?currentJobId=4364203279&keywords=Junior%20Software%20Engineer%20Ireland&location=Ireland")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
